$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ciudades")

# Update the "last updated" timestamp text in A1 (22:20 -> 23:20)
$ws.Range("A1").Value = "Datos actualizados a 2 de Abril de 2020 a las 23:20"

# Update Cataluña row (row 5) figures
$ws.Range("B5").Value = 23460
$ws.Range("C5").Value = 7849
$ws.Range("D5").Value = 13276
$ws.Range("E5").Value = 2335
